$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Rows 3 and 4: typed in directly (no explicit font color override)
$ws.Range("B3").Value = "Wouter Deprez"
$ws.Range("B4").Value = "Santa Cruz"
$ws.Range("A4").Value = "a32786c2-9043-4c1f-b81e-3f9a3d37e205"
$ws.Range("A3").Value = "c8d9ff35-6924-4911-b97a-f0d8e128796a"

# Rows 5-7: pasted in (carry an explicit black font color)
$ws.Range("B5").Value = "Xander De Rycke"
$ws.Range("A5").Value = "04a2d286-a0b7-40fe-a909-d6cf8c93c4a7"
$ws.Range("A6").Value = "20578d85-aac1-4866-aa21-ffd0dbaf5145"
$ws.Range("B6").Value = "Mat Bastard"
$ws.Range("B7").Value = "The Low Countries"
$ws.Range("A7").Value = "b22ed78a-4c25-4f3e-b233-1d76cdaf4ff5"

$ws.Range("A5:B7").Font.Color = 0

# Move the selection down past the newly added rows
$ws.Range("A8").Select()
